$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Bump the title row to the new sp_Blitz Check ID List version/date.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "sp_Blitz® Check ID List - v42 2015-09-07"

# ---------------------------------------------------------------------------
# 2. Append the seven new CheckID rows (158-164) below the existing data
#    (which currently ends at row 222 / CheckID 157).
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=223; CheckID=158; Priority=100; Group="Performance";  Finding="File growth set to 1MB";                 Url="http://BrentOzar.com/go/percentgrowth" },
    @{ Row=224; CheckID=159; Priority=1;   Group="Performance";  Finding="Memory Dangerously Low in NUMA Nodes";    Url="http://BrentOzar.com/go/max" },
    @{ Row=225; CheckID=160; Priority=100; Group="Performance";  Finding="Many Plans for One Query";                Url="http://BrentOzar.com/go/parameterization" },
    @{ Row=226; CheckID=161; Priority=100; Group="Performance";  Finding="High Number of Cached Plans";             Url="http://BrentOzar.com/go/planlimits" },
    @{ Row=227; CheckID=162; Priority=100; Group="Performance";  Finding="Poison Wait Detected: CMEMTHREAD & NUMA"; Url="http://BrentOzar.com/go/poison" },
    @{ Row=228; CheckID=163; Priority=10;  Group="Performance";  Finding="Query Store Disabled";                   Url="http://BrentOzar.com/go/querystore" },
    @{ Row=229; CheckID=164; Priority=20;  Group="Reliability";  Finding="Plan Guides Failing";                    Url="http://BrentOzar.com/go/guides" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.CheckID
    $ws.Cells.Item($row, 2).Value2 = $r.Priority
    $ws.Cells.Item($row, 3).Value2 = $r.Group
    $ws.Cells.Item($row, 4).Value2 = $r.Finding
    $ws.Cells.Item($row, 5).Value2 = $r.Url

    # Match the look of the surrounding table (column A centered numbers,
    # B:D plain text, E as a hyperlink-styled cell).
    $ws.Cells.Item($row, 1).Style = $ws.Range("A222").Style
    $ws.Cells.Item($row, 2).Style = $ws.Range("B222").Style
    $ws.Cells.Item($row, 3).Style = $ws.Range("C222").Style
    $ws.Cells.Item($row, 4).Style = $ws.Range("D222").Style
    $ws.Cells.Item($row, 5).Style = $ws.Range("E222").Style

    $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $r.Url) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. Keep the frozen header pane anchored near the top of the refreshed list
#    and reselect the first data row under the freeze.
# ---------------------------------------------------------------------------
$ws.Range("B5").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
